$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Updated" label is being localized to Spanish ("Actualizado").
# It appeared (wrapped) in column B on rows 2, 5, 6, 7 and 10 - update the
# text and drop the wrap-text formatting that used to highlight it.
$ws.Range("B2").Value = "Actualizado"
$ws.Range("B2").WrapText = $false

$ws.Range("B5").Value = "Actualizado"
$ws.Range("B5").WrapText = $false

$ws.Range("B6").Value = "Actualizado"
$ws.Range("B6").WrapText = $false

$ws.Range("B7").Value = "Actualizado"
$ws.Range("B7").WrapText = $false

$ws.Range("B10").Value = "Actualizado"
$ws.Range("B10").WrapText = $false

# Row 15 column B incorrectly duplicated column A's value ("Agatha"); fix it
# up to read "Actualizado" like the rest of column B, without the wrap.
$ws.Range("B15").Value = "Actualizado"
$ws.Range("B15").WrapText = $false

# The active selection moved from B10 to B15.
$ws.Range("B15").Select() | Out-Null
